# Automatic update of files.
# Rows 17-21 in the "Artfynd" sheet get re-ordered (re-sorted), carrying their
# Id (A), East/Ost (Q) and North/Nord (R) coordinates with them. Rows that
# represent records 111821923 and 111821924 also carry extra (mostly empty)
# attribute cells, among them K = "blomning".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 17-21 after the re-sort, in row order.
# Each entry: row number, Id (A), Ost (Q), Nord (R), whether it carries the
# extra J/K/L/N/AF attribute cells (K = "blomning").
$rowsData = @(
    @{ Row = 17; A = 111821924; Q = 550675.3931295178; R = 6681937.422269406; Extra = $true  },
    @{ Row = 18; A = 111821923; Q = 550701.1291094749; R = 6681909.496304798; Extra = $true  },
    @{ Row = 19; A = 111821927; Q = 550819.8901872271; R = 6681733.007140613; Extra = $false },
    @{ Row = 20; A = 111821928; Q = 550825.9503372401; R = 6681726.144349095; Extra = $false },
    @{ Row = 21; A = 111821926; Q = 550846.2444635418; R = 6681625.195240833; Extra = $false }
)

foreach ($rd in $rowsData) {
    $r = $rd.Row

    # Id, East (Ost) and North (Nord) coordinates.
    $ws.Range("A$r").Value = $rd.A
    $ws.Range("Q$r").Value = $rd.Q
    $ws.Range("R$r").Value = $rd.R

    if ($rd.Extra) {
        # These rows carry the (mostly empty) J/K/L/N/AF attribute cells;
        # only K has real content ("blomning" = Age/Stage "flowering").
        $ws.Range("J$r").Value = ""
        $ws.Range("K$r").Value = "blomning"
        $ws.Range("L$r").Value = ""
        $ws.Range("N$r").Value = ""
        $ws.Range("AF$r").Value = ""
    }
    else {
        # These rows must not carry the extra attribute cells.
        $ws.Range("J$r").ClearContents()
        $ws.Range("K$r").ClearContents()
        $ws.Range("L$r").ClearContents()
        $ws.Range("N$r").ClearContents()
        $ws.Range("AF$r").ClearContents()
    }
}
